# Generate Report for Handoff
# Updates the status of the two most-recently-handed-off files
# (92be5cd7-fa47-465e-9c9d-83665fb258a4 and be4fd1a7-9918-4402-8d16-3defea31a54c)
# from "Handed back: in sync with en-US" to "Ready for handoff", refreshes the
# related timestamps, and records a "stale handback version" error detail for
# each locale sheet.

$wb = $excel.ActiveWorkbook

$newStatus = "Ready for handoff"

# ---------------------------------------------------------------------------
# Overview sheet: rows 4 (92be5cd7...) and 5 (be4fd1a7...)
# Columns: E = zh-cn status, F = de-de status, G = Latest HO Xliff Generate Date
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E4").Value = $newStatus
$overview.Range("F4").Value = $newStatus
$overview.Range("G4").Value = "2016-09-03 22:29:08"

$overview.Range("E5").Value = $newStatus
$overview.Range("F5").Value = $newStatus
$overview.Range("G5").Value = "2016-09-03 22:29:08"

# ---------------------------------------------------------------------------
# zh-cn sheet: rows 4 and 5
# Columns: C = Status, H = Latest Handoff Datetime, P = Error Detail
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C4").Value = $newStatus
$zhcn.Range("H4").Value = "2016-09-03 22:28:59"
$zhcn.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f391cfa87a14ddf645061ec3c959d41b3e987ed1/e2e/92be5cd7-fa47-465e-9c9d-83665fb258a4.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/13c85e050f4eb7c4fb135e165f1fe0bdf56222f8/e2e/92be5cd7-fa47-465e-9c9d-83665fb258a4.md."

$zhcn.Range("C5").Value = $newStatus
$zhcn.Range("H5").Value = "2016-09-03 22:28:59"
$zhcn.Range("P5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f391cfa87a14ddf645061ec3c959d41b3e987ed1/e2e/be4fd1a7-9918-4402-8d16-3defea31a54c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/13c85e050f4eb7c4fb135e165f1fe0bdf56222f8/e2e/be4fd1a7-9918-4402-8d16-3defea31a54c.md."

# Widen the Error Detail column now that it holds a long message.
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------------
# de-de sheet: rows 4 and 5
# Columns: C = Status, H = Latest Handoff Datetime, P = Error Detail
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C4").Value = $newStatus
$dede.Range("H4").Value = "2016-09-03 22:29:08"
$dede.Range("P4").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f391cfa87a14ddf645061ec3c959d41b3e987ed1/e2e/92be5cd7-fa47-465e-9c9d-83665fb258a4.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/13c85e050f4eb7c4fb135e165f1fe0bdf56222f8/e2e/92be5cd7-fa47-465e-9c9d-83665fb258a4.md."

$dede.Range("C5").Value = $newStatus
$dede.Range("H5").Value = "2016-09-03 22:29:08"
$dede.Range("P5").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f391cfa87a14ddf645061ec3c959d41b3e987ed1/e2e/be4fd1a7-9918-4402-8d16-3defea31a54c.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/13c85e050f4eb7c4fb135e165f1fe0bdf56222f8/e2e/be4fd1a7-9918-4402-8d16-3defea31a54c.md."

# Widen the Error Detail column now that it holds a long message.
$dede.Columns.Item(16).ColumnWidth = 39.17
